$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (PORC_AVANCE), shifting it to G, and
# copying its style so the header keeps the same formatting.
$ws.Columns("F:F").Insert()

# Header
$ws.Range("F1").Value = "Monitor"

# Monitor names per row
$monitors = @(
    "Ariadna Chavez",
    "Verioska Butron",
    "Sofía Asto",
    "Alonso Chapoñan",
    "Valeria Palacios",
    "Ana Claudia Arana",
    "Ivanna Carrasco",
    "Betsy Fidel",
    "Mary Arcos",
    "Carito Hernandez ",
    "Banesa Perez",
    "Luis García",
    "Dorian Macha",
    "Mia Espejo"
)

for ($i = 0; $i -lt $monitors.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $monitors[$i]
}

# Match header style of the other header cells (bold, bordered, centered)
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1
